$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3.97069755430319
$ws.Range("D2").Value = 9.343553168875742
$ws.Range("E2").Value = 13.73693510877937
$ws.Range("F2").Value = 31.98497399215396
$ws.Range("G2").Value = 3.645236685970402
$ws.Range("I2").Value = 22.12485348301787
$ws.Range("J2").Value = 9.88353838016239
$ws.Range("K2").Value = 15.07115703360757
$ws.Range("N2").Value = 17.18494026279522
$ws.Range("O2").Value = 23.92024737470009
$ws.Range("C3").Value = 3.876193474139945
$ws.Range("D3").Value = 9.292105298099736
$ws.Range("E3").Value = 13.67717787110803
$ws.Range("F3").Value = 32.00132083471997
$ws.Range("G3").Value = 3.647436812006525
$ws.Range("I3").Value = 22.16263564821602
$ws.Range("J3").Value = 9.889969066870449
$ws.Range("K3").Value = 14.61000172502458
$ws.Range("N3").Value = 17.21683520241521
$ws.Range("O3").Value = 23.98198237641336
$ws.Range("C4").Value = 3.816303333258358
$ws.Range("D4").Value = 9.261869134761881
$ws.Range("E4").Value = 13.64320806552081
$ws.Range("F4").Value = 32.02028928498228
$ws.Range("G4").Value = 3.64885935151357
$ws.Range("I4").Value = 22.19171358719092
$ws.Range("J4").Value = 9.89554775896101
$ws.Range("K4").Value = 14.32080056799128
$ws.Range("N4").Value = 17.23824125943418
$ws.Range("O4").Value = 24.0260560528477
$ws.Range("C5").Value = 3.791448430726454
$ws.Range("D5").Value = 9.249897313255637
$ws.Range("E5").Value = 13.63005994388026
$ws.Range("F5").Value = 32.03026126785269
$ws.Range("G5").Value = 3.6494571236529
$ws.Range("I5").Value = 22.2050373706762
$ws.Range("J5").Value = 9.89823126846232
$ws.Range("K5").Value = 14.20160050145848
$ws.Range("N5").Value = 17.24742338215255
$ws.Range("O5").Value = 24.04556214361808
$ws.Range("C6").Value = 3.787294762007731
$ws.Range("D6").Value = 9.247930786143614
$ws.Range("E6").Value = 13.62791896592096
$ws.Range("F6").Value = 32.03205239602334
$ws.Range("G6").Value = 3.649557476696318
$ws.Range("I6").Value = 22.20733868534775
$ws.Range("J6").Value = 9.898701639016521
$ws.Range("K6").Value = 14.18173130293402
$ws.Range("N6").Value = 17.24897580899107
$ws.Range("O6").Value = 24.04889430970482
$ws.Range("C7").Value = 3.81596992322606
$ws.Range("D7").Value = 9.261706250723755
$ws.Range("E7").Value = 13.6430279186699
$ws.Range("F7").Value = 32.02041469803596
$ws.Range("G7").Value = 3.648867340012732
$ws.Range("I7").Value = 22.19188731301347
$ws.Range("J7").Value = 9.895582288842665
$ws.Range("K7").Value = 14.3191982149858
$ws.Range("N7").Value = 17.23836323341594
$ws.Range("O7").Value = 24.02631286762015
$ws.Range("C8").Value = 3.938510819139503
$ws.Range("D8").Value = 9.325538917349926
$ws.Range("E8").Value = 13.71577257724662
$ws.Range("F8").Value = 31.98875484498594
$ws.Range("G8").Value = 3.645980452745189
$ws.Range("I8").Value = 22.13665813125263
$ws.Range("J8").Value = 9.885417488370107
$ws.Range("K8").Value = 14.91350611327996
$ws.Range("N8").Value = 17.19555981886677
$ws.Range("O8").Value = 23.94025031472835
$ws.Range("C9").Value = 4.163228654192226
$ws.Range("D9").Value = 9.461015574272993
$ws.Range("E9").Value = 13.87951510752249
$ws.Range("F9").Value = 31.99764383734275
$ws.Range("G9").Value = 3.640885166627954
$ws.Range("I9").Value = 22.07516856557783
$ws.Range("J9").Value = 9.878406419917713
$ws.Range("K9").Value = 16.02400689956949
$ws.Range("N9").Value = 17.12605177586768
$ws.Range("O9").Value = 23.82064633501019
$ws.Range("C10").Value = 4.317862430201219
$ws.Range("D10").Value = 9.566220662833834
$ws.Range("E10").Value = 14.01194687810959
$ws.Range("F10").Value = 32.04750280776449
$ws.Range("G10").Value = 3.637482938653894
$ws.Range("I10").Value = 22.05871987144928
$ws.Range("J10").Value = 9.88111009754468
$ws.Range("K10").Value = 16.79806282003841
$ws.Range("N10").Value = 17.08373991632443
$ws.Range("O10").Value = 23.76302774906321
$ws.Range("C11").Value = 4.385758617986729
$ws.Range("D11").Value = 9.615175279199088
$ws.Range("E11").Value = 14.07465682068189
$ws.Range("F11").Value = 32.07957564892024
$ws.Range("G11").Value = 3.636008494407183
$ws.Range("I11").Value = 22.05750214222378
$ws.Range("J11").Value = 9.884038242897534
$ws.Range("K11").Value = 17.13961854440843
$ws.Range("N11").Value = 17.0663838598383
$ws.Range("O11").Value = 23.74343684740311
$ws.Range("C12").Value = 4.411103958922822
$ws.Range("D12").Value = 9.633858749612394
$ws.Range("E12").Value = 14.09874263118752
$ws.Range("F12").Value = 32.09306720422364
$ws.Range("G12").Value = 3.635460633688095
$ws.Range("I12").Value = 22.05794297051788
$ws.Range("J12").Value = 9.885390380072367
$ws.Range("K12").Value = 17.26732926442858
$ws.Range("N12").Value = 17.06008293142805
$ws.Range("O12").Value = 23.73697371384501
$ws.Range("C13").Value = 4.405661843211455
$ws.Range("D13").Value = 9.629828659381392
$ws.Range("E13").Value = 14.09354050370092
$ws.Range("F13").Value = 32.09010176512471
$ws.Range("G13").Value = 3.635578160071965
$ws.Range("I13").Value = 22.05780790588769
$ws.Range("J13").Value = 9.885088367406919
$ws.Range("K13").Value = 17.23989854737675
$ws.Range("N13").Value = 17.06142788667461
$ws.Range("O13").Value = 23.73832311713939
$ws.Range("C14").Value = 4.387851192203984
$ws.Range("D14").Value = 9.616709535937828
$ws.Range("E14").Value = 14.07663167441578
$ws.Range("F14").Value = 32.0806586560136
$ws.Range("G14").Value = 3.635963211859666
$ws.Range("I14").Value = 22.05752032923515
$ws.Range("J14").Value = 9.8841446132878
$ws.Range("K14").Value = 17.1501586162141
$ws.Range("N14").Value = 17.06586004220836
$ws.Range("O14").Value = 23.74288595295583
$ws.Range("C15").Value = 4.376893678275164
$ws.Range("D15").Value = 9.608692265311534
$ws.Range("E15").Value = 14.06631818920667
$ws.Range("F15").Value = 32.07504964980878
$ws.Range("G15").Value = 3.636200430335521
$ws.Range("I15").Value = 22.05746166143404
$ws.Range("J15").Value = 9.883598194284986
$ws.Range("K15").Value = 17.09497511888932
$ws.Range("N15").Value = 17.06861019610069
$ws.Range("O15").Value = 23.74580535506636
$ws.Range("C16").Value = 4.31337482079196
$ws.Range("D16").Value = 9.563042383885989
$ws.Range("E16").Value = 14.00789693082899
$ws.Range("F16").Value = 32.04559542988321
$ws.Range("G16").Value = 3.637580766381139
$ws.Range("I16").Value = 22.05892562602319
$ws.Range("J16").Value = 9.880952835991879
$ws.Range("K16").Value = 16.77551908978731
$ws.Range("N16").Value = 17.08491219166768
$ws.Range("O16").Value = 23.76444158877264
$ws.Range("C17").Value = 4.273771185108491
$ws.Range("D17").Value = 9.53530982342239
$ws.Range("E17").Value = 13.97267801974764
$ws.Range("F17").Value = 32.02992913786434
$ws.Range("G17").Value = 3.63844627984152
$ws.Range("I17").Value = 22.06142929103657
$ws.Range("J17").Value = 9.879764360603922
$ws.Range("K17").Value = 16.5767551190346
$ws.Range("N17").Value = 17.0953970451445
$ws.Range("O17").Value = 23.77757256245817
$ws.Range("C18").Value = 4.250762581048889
$ws.Range("D18").Value = 9.519462724419382
$ws.Range("E18").Value = 13.95265413113731
$ws.Range("F18").Value = 32.02180264992898
$ws.Range("G18").Value = 3.638950997749052
$ws.Range("I18").Value = 22.06345899129835
$ws.Range("J18").Value = 9.879240688257964
$ws.Range("K18").Value = 16.46144257242087
$ws.Range("N18").Value = 17.10160575879337
$ws.Range("O18").Value = 23.78574808260127
$ws.Range("C19").Value = 4.242933254274874
$ws.Range("D19").Value = 9.514115388845749
$ws.Range("E19").Value = 13.9459148898286
$ws.Range("F19").Value = 32.01920315273244
$ws.Range("G19").Value = 3.639123072839565
$ws.Range("I19").Value = 22.06424744308572
$ws.Range("J19").Value = 9.879090870954087
$ws.Range("K19").Value = 16.42223343901523
$ws.Range("N19").Value = 17.1037385305407
$ws.Range("O19").Value = 23.78862304594793
$ws.Range("C20").Value = 4.278010928476875
$ws.Range("D20").Value = 9.538251332657138
$ws.Range("E20").Value = 13.97640312010038
$ws.Range("F20").Value = 32.03150534044919
$ws.Range("G20").Value = 3.638353430982458
$ws.Range("I20").Value = 22.0611017341407
$ws.Range("J20").Value = 9.8798743320603
$ws.Range("K20").Value = 16.59801704838593
$ws.Range("N20").Value = 17.09426248597783
$ws.Range("O20").Value = 23.77611024947259
$ws.Range("C21").Value = 4.393092632722404
$ws.Range("D21").Value = 9.620559092383019
$ws.Range("E21").Value = 14.08158913736846
$ws.Range("F21").Value = 32.08339583039923
$ws.Range("G21").Value = 3.635849828833933
$ws.Range("I21").Value = 22.05758031356027
$ws.Range("J21").Value = 9.884415221106172
$ws.Range("K21").Value = 17.17656243566995
$ws.Range("N21").Value = 17.06455084889626
$ws.Range("O21").Value = 23.74151977916961
$ws.Range("C22").Value = 4.466169499478869
$ws.Range("D22").Value = 9.67519377312049
$ws.Range("E22").Value = 14.15230206663046
$ws.Range("F22").Value = 32.12515326552224
$ws.Range("G22").Value = 3.6342746370816
$ws.Range("I22").Value = 22.06053615023656
$ws.Range("J22").Value = 9.888800654098178
$ws.Range("K22").Value = 17.54513742248717
$ws.Range("N22").Value = 17.04671441823849
$ws.Range("O22").Value = 23.72448367972873
$ws.Range("C23").Value = 4.427366572871212
$ws.Range("D23").Value = 9.64596128424026
$ws.Range("E23").Value = 14.11438652628399
$ws.Range("F23").Value = 32.10215061741205
$ws.Range("G23").Value = 3.635109777390897
$ws.Range("I23").Value = 22.05847734948518
$ws.Range("J23").Value = 9.886330675949987
$ws.Range("K23").Value = 17.34932793605365
$ws.Range("N23").Value = 17.05608951860571
$ws.Range("O23").Value = 23.73306539488724
$ws.Range("C24").Value = 4.276094886849827
$ws.Range("D24").Value = 9.536921174145176
$ws.Range("E24").Value = 13.97471830353454
$ws.Range("F24").Value = 32.03078999692396
$ws.Range("G24").Value = 3.638395385760848
$ws.Range("I24").Value = 22.06124798389266
$ws.Range("J24").Value = 9.879824116767141
$ws.Range("K24").Value = 16.58840776255596
$ws.Range("N24").Value = 17.09477485683623
$ws.Range("O24").Value = 23.77676941036694
$ws.Range("C25").Value = 4.10420978749578
$ws.Range("D25").Value = 9.423324146883166
$ws.Range("E25").Value = 13.83303247962456
$ws.Range("F25").Value = 31.98762876962223
$ws.Range("G25").Value = 3.642203378347935
$ws.Range("I25").Value = 22.08676936116854
$ws.Range("J25").Value = 9.878921561004923
$ws.Range("K25").Value = 15.73038091970231
$ws.Range("N25").Value = 17.14331501829181
$ws.Range("O25").Value = 23.82064633501019
